$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3126816
$ws.Range("J17").Value = 3126816
$ws.Range("L17").Value = 9380448
$ws.Range("N17").Value = -9380784

$ws.Range("H87").Value = 15818.066
$ws.Range("J87").Value = 15818.066
$ws.Range("L87").Value = 15818.066
$ws.Range("N87").Value = -18314.066

$ws.Range("H90").Value = 15818.066
$ws.Range("J90").Value = 15818.066
$ws.Range("L90").Value = 47454.198
$ws.Range("N90").Value = -59934.198

$ws.Range("H98").Value = 1492.091
$ws.Range("I98").Value = 1485.8182
$ws.Range("K98").Value = 1485.8182
$ws.Range("M98").Value = 12.18180000000007

$ws.Range("H112").Value = 1096.7709
$ws.Range("I112").Value = 792.3333
$ws.Range("J112").Value = 1117.0667
$ws.Range("K112").Value = 2376.9999
$ws.Range("L112").Value = 3351.2001
$ws.Range("M112").Value = -1268.9999
$ws.Range("N112").Value = -5567.2001

$ws.Range("H122").Value = 1492.091
$ws.Range("I122").Value = 1485.8182
$ws.Range("K122").Value = 4457.4546
$ws.Range("M122").Value = -2007.4546

$ws.Range("H129").Value = 805.76746
$ws.Range("I129").Value = 485
$ws.Range("K129").Value = 1455
$ws.Range("M129").Value = 3545

$ws.Range("H137").Value = 728.2
$ws.Range("I137").Value = 732.8333
$ws.Range("J137").Value = 686.5
$ws.Range("K137").Value = 2198.4999
$ws.Range("L137").Value = 2059.5
$ws.Range("M137").Value = 351.5001000000002
$ws.Range("N137").Value = -7159.5

$ws.Range("H138").Value = 3416.46
$ws.Range("I138").Value = 1902.3715
$ws.Range("J138").Value = 4231.7383
$ws.Range("K138").Value = 5707.1145
$ws.Range("L138").Value = 12695.2149
$ws.Range("M138").Value = -567.1144999999997
$ws.Range("N138").Value = -22975.2149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6746.174
$ws.Range("I74").Value = 10570.667
$ws.Range("K74").Value = 10570.667
$ws.Range("M74").Value = -9696.666999999999

$ws.Range("H77").Value = 6746.174
$ws.Range("I77").Value = 10570.667
$ws.Range("K77").Value = 52853.335
$ws.Range("M77").Value = -48485.335

$ws.Range("H132").Value = 2178
$ws.Range("I132").Value = 1561.4615
$ws.Range("J132").Value = 2845.9167
$ws.Range("K132").Value = 4684.3845
$ws.Range("L132").Value = 8537.750100000001
$ws.Range("M132").Value = -2154.3845
$ws.Range("N132").Value = -13597.7501

$ws.Range("H135").Value = 45666.668
$ws.Range("J135").Value = 45666.668
$ws.Range("L135").Value = 45666.668
$ws.Range("N135").Value = -55806.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2807.875
$ws.Range("I80").Value = 1117.4
$ws.Range("J80").Value = 3252.7368
$ws.Range("K80").Value = 1117.4
$ws.Range("L80").Value = 3252.7368
$ws.Range("M80").Value = -119.4000000000001
$ws.Range("N80").Value = -5248.736800000001

$ws.Range("H83").Value = 2807.875
$ws.Range("I83").Value = 1117.4
$ws.Range("J83").Value = 3252.7368
$ws.Range("K83").Value = 5587
$ws.Range("L83").Value = 16263.684
$ws.Range("M83").Value = -595
$ws.Range("N83").Value = -26247.684

$ws.Range("H134").Value = 1903.7916
$ws.Range("I134").Value = 1741.6842
$ws.Range("K134").Value = 5225.0526
$ws.Range("M134").Value = -2690.0526

$ws.Range("H135").Value = 40780
$ws.Range("J135").Value = 40780
$ws.Range("L135").Value = 40780
$ws.Range("N135").Value = -50920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1725.7
$ws.Range("I31").Value = 917.9853000000001
$ws.Range("J31").Value = 3442.0938
$ws.Range("K31").Value = 917.9853000000001
$ws.Range("L31").Value = 3442.0938
$ws.Range("M31").Value = -622.9853000000001
$ws.Range("N31").Value = -4032.0938

$ws.Range("H34").Value = 1725.7
$ws.Range("I34").Value = 917.9853000000001
$ws.Range("J34").Value = 3442.0938
$ws.Range("K34").Value = 917.9853000000001
$ws.Range("L34").Value = 3442.0938
$ws.Range("M34").Value = -715.9853000000001
$ws.Range("N34").Value = -3846.0938

$ws.Range("H58").Value = 2466.037
$ws.Range("J58").Value = 2653.1538
$ws.Range("L58").Value = 2653.1538
$ws.Range("N58").Value = -3059.1538

$ws.Range("H122").Value = 3218.8572
$ws.Range("I122").Value = 2183.6
$ws.Range("J122").Value = 5807
$ws.Range("K122").Value = 6550.799999999999
$ws.Range("L122").Value = 17421
$ws.Range("M122").Value = -4100.799999999999
$ws.Range("N122").Value = -22321

$ws.Range("H134").Value = 1790.258
$ws.Range("I134").Value = 1891.6923
$ws.Range("J134").Value = 1262.8
$ws.Range("K134").Value = 5675.0769
$ws.Range("L134").Value = 3788.4
$ws.Range("M134").Value = -3140.0769
$ws.Range("N134").Value = -8858.4

$ws.Range("H136").Value = 2466.037
$ws.Range("J136").Value = 2653.1538
$ws.Range("L136").Value = 7959.4614
$ws.Range("N136").Value = -13059.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 97.26667
$ws.Range("I23").Value = 55
$ws.Range("J23").Value = 100.28571
$ws.Range("K23").Value = 165
$ws.Range("L23").Value = 300.85713
$ws.Range("M23").Value = 70
$ws.Range("N23").Value = -770.85713

$ws.Range("H92").Value = 988
$ws.Range("J92").Value = 925
$ws.Range("L92").Value = 2775
$ws.Range("N92").Value = -5271

$ws.Range("H104").Value = 2938
$ws.Range("I104").Value = 1000
$ws.Range("J104").Value = 3422.5
$ws.Range("K104").Value = 3000
$ws.Range("L104").Value = 10267.5
$ws.Range("M104").Value = -379
$ws.Range("N104").Value = -15509.5

$ws.Range("H105").Value = 8466.666999999999
$ws.Range("J105").Value = 8466.666999999999
$ws.Range("L105").Value = 25400.001
$ws.Range("N105").Value = -30642.001

$ws.Range("H106").Value = 5244.4443
$ws.Range("J106").Value = 5244.4443
$ws.Range("L106").Value = 15733.3329
$ws.Range("N106").Value = -17625.3329

$ws.Range("H107").Value = 801933.2
$ws.Range("I107").Value = 2710.6
$ws.Range("J107").Value = 1134942.6
$ws.Range("K107").Value = 8131.799999999999
$ws.Range("L107").Value = 3404827.8
$ws.Range("M107").Value = -6211.799999999999
$ws.Range("N107").Value = -3408667.8

$ws.Range("H108").Value = 5730
$ws.Range("I108").Value = 1990
$ws.Range("J108").Value = 7600
$ws.Range("K108").Value = 5970
$ws.Range("L108").Value = 22800
$ws.Range("M108").Value = -3090
$ws.Range("N108").Value = -28560

$ws.Range("H109").Value = 5902.5713
$ws.Range("J109").Value = 6374
$ws.Range("L109").Value = 19122
$ws.Range("N109").Value = -21202

$ws.Range("H113").Value = 1247.4865
$ws.Range("I113").Value = 1238.72
$ws.Range("J113").Value = 1265.75
$ws.Range("K113").Value = 3716.16
$ws.Range("L113").Value = 3797.25
$ws.Range("M113").Value = -1546.16
$ws.Range("N113").Value = -8137.25

$ws.Range("H122").Value = 1196.1428
$ws.Range("I122").Value = 465.13043
$ws.Range("J122").Value = 4558.8
$ws.Range("K122").Value = 4186.17387
$ws.Range("L122").Value = 41029.2
$ws.Range("M122").Value = -1736.17387
$ws.Range("N122").Value = -45929.2

$ws.Range("H123").Value = 1416.6666
$ws.Range("I123").Value = 800
$ws.Range("J123").Value = 1493.75
$ws.Range("K123").Value = 2400
$ws.Range("L123").Value = 4481.25
$ws.Range("M123").Value = 50
$ws.Range("N123").Value = -9381.25

$ws.Range("H124").Value = 1826.4286
$ws.Range("I124").Value = 600
$ws.Range("J124").Value = 1920.7693
$ws.Range("K124").Value = 1800
$ws.Range("L124").Value = 5762.3079
$ws.Range("M124").Value = 3110
$ws.Range("N124").Value = -15582.3079

$ws.Range("H125").Value = 1320
$ws.Range("I125").Value = 666.6667
$ws.Range("J125").Value = 1600
$ws.Range("K125").Value = 2000.0001
$ws.Range("L125").Value = 4800
$ws.Range("M125").Value = 2919.9999
$ws.Range("N125").Value = -14640

$ws.Range("H131").Value = 854.33
$ws.Range("J131").Value = 857.9091
$ws.Range("L131").Value = 2573.7273
$ws.Range("N131").Value = -12653.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 67071.42999999999
$ws.Range("J139").Value = 67071.42999999999
$ws.Range("L139").Value = 67071.42999999999
$ws.Range("N139").Value = -77351.42999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 832.84
$ws.Range("I16").Value = 771.75
$ws.Range("J16").Value = 889.2308
$ws.Range("K16").Value = 771.75
$ws.Range("L16").Value = 889.2308
$ws.Range("M16").Value = -601.75
$ws.Range("N16").Value = -1229.2308

$ws.Range("H55").Value = 364.86365
$ws.Range("I55").Value = 207.8
$ws.Range("J55").Value = 495.75
$ws.Range("K55").Value = 207.8
$ws.Range("L55").Value = 495.75
$ws.Range("M55").Value = -34.80000000000001
$ws.Range("N55").Value = -841.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 49214.5
$ws.Range("J138").Value = 49214.5
$ws.Range("L138").Value = 49214.5
$ws.Range("N138").Value = -59494.5
